# Fix conversion error (off by 1000) in land/BLAPE BAU LULUCF Anthro Pollutant Emis.xlsx
# Row 2 ("CO2 (g)") of the BLAPE sheet multiplies EEA (MtCO2e) values by 1e9 to get grams.
# That was off by a factor of 1000 (should be *1e12), so update the formulas for every
# cell in row 2 that is NOT already correct (cell C2 already used the right multiplier).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLAPE")

$cols = @("B","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

foreach ($col in $cols) {
    $cell = $ws.Range($col + "2")
    $oldFormula = $cell.Formula
    $newFormula = $oldFormula -replace [regex]::Escape("*1000000000"), "*1000000000000"
    $cell.Formula = $newFormula
}

$excel.CalculateFullRebuild()
